$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet: bump the "Last Updated" timestamp by one minute.
# ---------------------------------------------------------------------------
$wsMeta = $wb.Sheets.Item("Metadata")
$wsMeta.Range("A2").Value = "05 Nov 2025, 10:16 AM"

# ---------------------------------------------------------------------------
# 2. Industry Analysis sheet: refresh the "1 Year" column (F2:F76).
# ---------------------------------------------------------------------------
$wsIndustry = $wb.Sheets.Item("Industry Analysis")

$oneYearValues = @(
    21.3,
    -4.3927,
    35.9445,
    -51.0482,
    57.2275,
    -9.6407,
    -6.1449,
    36.9733,
    -4.7026,
    46.5317,
    -2.102,
    17.4681,
    -33.0245,
    1.0205,
    2.0426,
    -16.2411,
    7.4627,
    -25.798,
    47.7485,
    19.5587,
    76.5603,
    -54.2675,
    -0.8811,
    4.8518,
    3.6831,
    -34.0874,
    -11.9893,
    -12.994,
    25.5415,
    56.5088,
    2.0908,
    -4.7193,
    22.8807,
    5.3359,
    -5.1995,
    -5.6238,
    -22.595,
    10.8405,
    -7.5963,
    -4.552,
    22.3098,
    14.0694,
    -9.6066,
    27.639,
    -6.3484,
    -40.5302,
    -29.7988,
    -24.0791,
    -49.1803,
    -51.6023,
    -34.4756,
    -11.5478,
    -2.3796,
    -15.4382,
    -27.6987,
    -27.1559,
    -2.1585,
    -23.0964,
    -13.3217,
    -8.1496,
    -16.0695,
    -12.5465,
    47.7264,
    -42.4232,
    11.3291,
    14.3746,
    32.6702,
    -17.0097,
    -13.5162,
    11.4259,
    2.6754,
    -11.1574,
    -13.2502,
    24.7078,
    53.3554
)

for ($i = 0; $i -lt $oneYearValues.Length; $i++) {
    $rowNum = 2 + $i
    $wsIndustry.Cells.Item($rowNum, 6).Value = $oneYearValues[$i]
}

# ---------------------------------------------------------------------------
# 3. Stock List sheet: a new row (CAPTRU-RE1) is inserted at row 2, pushing
#    every other row down by one and dropping the final row (TRAVELFOOD)
#    off the bottom of the A1:H76 range.
# ---------------------------------------------------------------------------
$wsStock = $wb.Sheets.Item("Stock List")

$wsStock.Rows(2).Insert()
$wsStock.Range("A2:H2").ClearFormats()

$wsStock.Cells.Item(2, 1).Value = [char]0x1F4CB
$wsStock.Cells.Item(2, 2).Value = "CAPTRU-RE1"
$wsStock.Cells.Item(2, 3).Value = "CAPTRU-RE1"
$wsStock.Cells.Item(2, 4).Value = 5.67
$wsStock.Cells.Item(2, 5).Value = -11.9565
$wsStock.Cells.Item(2, 6).Value = "N/A"
$wsStock.Cells.Item(2, 7).Value = "N/A"
$wsStock.Cells.Item(2, 8).Value = 0

# Remove the now-duplicated trailing row that the insert pushed past row 76.
$wsStock.Rows(77).Delete()
